# Update Xcl1-Xcr1 NATMI LR-pairs sheet with new TPM-derived values.
# A new target cluster ("MuSCs") is introduced, so the original 2 sending
# clusters (FAPs, Inflammatory-Mac) x 3 target clusters (ECs,
# Inflammatory-Mac, Resolving-Mac) = 6 data rows grows to 2 x 4 = 8 rows
# (ECs, Inflammatory-Mac, MuSCs, Resolving-Mac), and every numeric column
# is refreshed with the new TPM computation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("FAPs", "Xcl1", "Xcr1", "ECs", 1, 0.3333333333333333, 0.2175426666666667, 0.652628, 0.2420741369384981, 0.2420741369384981, 1, 0.3333333333333333, 0.044302, 0.132906, 0.02358675946309446, 0.02371390862331898, 0.009637575218666667, 0.086738176968, 0.005709744440204547, 0.00574052396342835),
    @("FAPs", "Xcl1", "Xcr1", "Inflammatory-Mac", 1, 0.3333333333333333, 0.2175426666666667, 0.652628, 0.2420741369384981, 0.2420741369384981, 3, 1, 1.334208333333333, 4.002625, 0.7103438000990813, 0.7141730509037375, 0.2902472387222222, 2.6122251485, 0.1719558623385981, 0.1728828249222563),
    @("FAPs", "Xcl1", "Xcr1", "MuSCs", 1, 0.3333333333333333, 0.2175426666666667, 0.652628, 0.2420741369384981, 0.2420741369384981, 1, 0.5, 0.0302125, 0.060425, 0.01608539050784934, 0.01078140135557499, 0.006572507816666667, 0.0394350469, 0.003893857024506339, 0.002609898428138369),
    @("FAPs", "Xcl1", "Xcr1", "Resolving-Mac", 1, 0.3333333333333333, 0.2175426666666667, 0.652628, 0.2420741369384981, 0.2420741369384981, 2, 0.6666666666666666, 0.4695343333333333, 1.408603, 0.2499840499299749, 0.2513316391173685, 0.1021437509648889, 0.9192937586839999, 0.06051467313518911, 0.06084088962467506),
    @("Inflammatory-Mac", "Xcl1", "Xcr1", "ECs", 2, 0.6666666666666666, 0.6811186666666668, 2.043356, 0.7579258630615019, 0.7579258630615019, 1, 0.3333333333333333, 0.044302, 0.132906, 0.02358675946309446, 0.02371390862331898, 0.03017491917066667, 0.271574272536, 0.01787701502288992, 0.01797338465989063),
    @("Inflammatory-Mac", "Xcl1", "Xcr1", "Inflammatory-Mac", 2, 0.6666666666666666, 0.6811186666666668, 2.043356, 0.7579258630615019, 0.7579258630615019, 3, 1, 1.334208333333333, 4.002625, 0.7103438000990813, 0.7141730509037375, 0.9087542010555557, 8.178787809500001, 0.5383879377604831, 0.5412902259814811),
    @("Inflammatory-Mac", "Xcl1", "Xcr1", "MuSCs", 2, 0.6666666666666666, 0.6811186666666668, 2.043356, 0.7579258630615019, 0.7579258630615019, 1, 0.5, 0.0302125, 0.060425, 0.01608539050784934, 0.01078140135557499, 0.02057829771666667, 0.1234697863, 0.012191533483343, 0.00817150292743662),
    @("Inflammatory-Mac", "Xcl1", "Xcr1", "Resolving-Mac", 2, 0.6666666666666666, 0.6811186666666668, 2.043356, 0.7579258630615019, 0.7579258630615019, 2, 0.6666666666666666, 0.4695343333333333, 1.408603, 0.2499840499299749, 0.2513316391173685, 0.3198085990742222, 2.878277391668, 0.1894693767947858, 0.1904907494926935)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $values = $rows[$i]
    for ($c = 1; $c -le $values.Count; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $values[$c - 1]
    }
}
